$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in I1 and J1, matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill columns I and J for data rows 2-31:
#   I = 1 (constant)
#   J = same value as column H on that row
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
